$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string as TEXT (not auto-converted to a number).
# Route the literal through a formula ( ="..." ) so Excel keeps it typed as text,
# then Copy + PasteSpecial(values-only) to flatten the formula to a literal value
# without touching the cell number format / style.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $escaped = $val.Replace('"', '""')
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue "D2" "64.470.32"
$ws.Range("E2").Value = "  +1.27%  "
Set-TextValue "D3" "3.081.02"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "558.56"
$ws.Range("E5").Value = "  +1.13%  "
Set-TextValue "D6" "143.74"
$ws.Range("E6").Value = "  +3.64%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.26%  "
Set-TextValue "D8" "3.078.35"
$ws.Range("E8").Value = "  -0.11%  "
Set-TextValue "D9" "0.498"
$ws.Range("E9").Value = "  -0.54%  "
Set-TextValue "D10" "6.36"
$ws.Range("E10").Value = "  +2.26%  "
Set-TextValue "D11" "0.151"
$ws.Range("E11").Value = "  -0.42%  "
Set-TextValue "D12" "0.471"
$ws.Range("E12").Value = "  +3.38%  "
Set-TextValue "D13" "0.0000228"
$ws.Range("E13").Value = "  +0.55%  "
Set-TextValue "D14" "35.15"
$ws.Range("E14").Value = "  -0.18%  "
Set-TextValue "D15" "3.597.42"
$ws.Range("E15").Value = "  +0.86%  "
Set-TextValue "D16" "64.609.70"
$ws.Range("E16").Value = "  +1.39%  "
Set-TextValue "D17" "3.084.82"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  +0.59%  "
Set-TextValue "D19" "6.76"
$ws.Range("E19").Value = "  -0.03%  "
Set-TextValue "D20" "476.72"
$ws.Range("E20").Value = "  -2.23%  "
Set-TextValue "D21" "13.80"
$ws.Range("E21").Value = "  +1.92%  "
Set-TextValue "D22" "0.682"
$ws.Range("E22").Value = "  -0.73%  "
Set-TextValue "D23" "7.52"
$ws.Range("E23").Value = "  +4.00%  "
Set-TextValue "D24" "13.45"
$ws.Range("E24").Value = "  +8.57%  "
Set-TextValue "D25" "80.87"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  -0.13%  "
Set-TextValue "D27" "2.78"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  +2.01%  "
Set-TextValue "D29" "2.05"
$ws.Range("E29").Value = "  +2.98%  "
Set-TextValue "D30" "0.999"
$ws.Range("E30").Value = "  -0.15%  "
Set-TextValue "D31" "25.99"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("E32").Value = "  +0.79%  "
Set-TextValue "D33" "2.49"
$ws.Range("E33").Value = "  +2.78%  "
Set-TextValue "D34" "5.58"
$ws.Range("E34").Value = "  -4.57%  "
Set-TextValue "D35" "6.12"
$ws.Range("E35").Value = "  +2.08%  "
Set-TextValue "D36" "54.79"
$ws.Range("E36").Value = "  -1.62%  "
Set-TextValue "D37" "467.68"
$ws.Range("E37").Value = "  -0.47%  "
Set-TextValue "D38" "0.0831"
$ws.Range("E38").Value = "  +1.57%  "
Set-TextValue "D39" "0.0407"
$ws.Range("E39").Value = "  +2.26%  "
Set-TextValue "D40" "2.95"
$ws.Range("E40").Value = "  +16.23%  "
Set-TextValue "D41" "2.965.93"
$ws.Range("E41").Value = "  -7.19%  "
Set-TextValue "D42" "8.23"
$ws.Range("E42").Value = "  -0.24%  "
Set-TextValue "D43" "0.114"
$ws.Range("E43").Value = "  -5.37%  "
Set-TextValue "D44" "28.11"
$ws.Range("E44").Value = "  +0.89%  "
Set-TextValue "D45" "0.258"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  +4.77%  "
$ws.Range("E48").Value = "  +1.68%  "
Set-TextValue "D49" "0.0₃0522"
$ws.Range("E49").Value = "  -0.73%  "
Set-TextValue "D50" "117.64"
$ws.Range("E50").Value = "  +0.92%  "
Set-TextValue "D51" "2.06"
$ws.Range("E51").Value = "  -1.02%  "
